$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.041283512961458
$ws.Range("D2").Value = 1.043523996761665
$ws.Range("E2").Value = 1.039520552964206
$ws.Range("F2").Value = 1.04005079143382
$ws.Range("I2").Value = 1.039689566532417
$ws.Range("J2").Value = 1.046365555210746
$ws.Range("K2").Value = 1.046297620559727
$ws.Range("L2").Value = 1.042305506525003
$ws.Range("M2").Value = 1.042834238978154

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.042611367692881
$ws.Range("D3").Value = 1.044557176498724
$ws.Range("E3").Value = 1.040663462255592
$ws.Range("F3").Value = 1.042007747970503
$ws.Range("I3").Value = 1.040095466089856
$ws.Range("J3").Value = 1.04733749022168
$ws.Range("K3").Value = 1.047141272094567
$ws.Range("L3").Value = 1.043257774557176
$ws.Range("M3").Value = 1.044598523683211

$ws.Range("B4").Value = 1.019999999999999
$ws.Range("C4").Value = 1.043469050889066
$ws.Range("D4").Value = 1.045224334847763
$ws.Range("E4").Value = 1.04140189179445
$ws.Range("F4").Value = 1.043272232234967
$ws.Range("I4").Value = 1.04035613510211
$ws.Range("J4").Value = 1.047964418026637
$ws.Range("K4").Value = 1.04768519380586
$ws.Range("L4").Value = 1.043872271983033
$ws.Range("M4").Value = 1.045737944252977

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043829261698173
$ws.Range("D5").Value = 1.045504482618464
$ws.Range("E5").Value = 1.041712066810589
$ws.Range("F5").Value = 1.043803405673268
$ws.Range("I5").Value = 1.040465249799803
$ws.Range("J5").Value = 1.048227509637835
$ws.Range("K5").Value = 1.047913389285166
$ws.Range("L5").Value = 1.044130208137158
$ws.Range("M5").Value = 1.046216445596791

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.043889721802199
$ws.Range("D6").Value = 1.04555150169296
$ws.Range("E6").Value = 1.041764131407691
$ws.Range("F6").Value = 1.043892568072879
$ws.Range("I6").Value = 1.040483543110529
$ws.Range("J6").Value = 1.048271656527891
$ws.Range("K6").Value = 1.047951676915581
$ws.Range("L6").Value = 1.044173493510831
$ws.Range("M6").Value = 1.046296758488863

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.043473865446935
$ws.Range("D7").Value = 1.045228079470074
$ws.Range("E7").Value = 1.041406037386923
$ws.Range("F7").Value = 1.04327933141609
$ws.Range("I7").Value = 1.040357594944178
$ws.Range("J7").Value = 1.047967935305718
$ws.Range("K7").Value = 1.04768824480359
$ws.Range("L7").Value = 1.043875720099084
$ws.Range("M7").Value = 1.045744340003313

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041732587914414
$ws.Range("D8").Value = 1.04387345262166
$ws.Range("E8").Value = 1.039907037066105
$ws.Range("F8").Value = 1.040712535034916
$ws.Range("I8").Value = 1.039827152543523
$ws.Range("J8").Value = 1.04669443799557
$ws.Range("K8").Value = 1.046583148454251
$ws.Range("L8").Value = 1.042627681304075
$ws.Range("M8").Value = 1.043430948557377

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.038652241919602
$ws.Range("D9").Value = 1.041475678783609
$ws.Range("E9").Value = 1.037256896222565
$ws.Range("F9").Value = 1.036175090221767
$ws.Range("I9").Value = 1.038877214315632
$ws.Range("J9").Value = 1.04443499111193
$ws.Range("K9").Value = 1.044620497051295
$ws.Range("L9").Value = 1.040415379287334
$ws.Range("M9").Value = 1.039337096938423

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036590202625463
$ws.Range("D10").Value = 1.039869661844941
$ws.Range("E10").Value = 1.035483974619461
$ws.Range("F10").Value = 1.033139512605955
$ws.Range("I10").Value = 1.038233535140552
$ws.Range("J10").Value = 1.042918035364776
$ws.Range("K10").Value = 1.043301491309166
$ws.Range("L10").Value = 1.038931417905413
$ws.Range("M10").Value = 1.036595362593583

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035695214086733
$ws.Range("D11").Value = 1.039172397650024
$ws.Range("E11").Value = 1.034714750064464
$ws.Range("F11").Value = 1.031822351539389
$ws.Range("I11").Value = 1.037952318176816
$ws.Range("J11").Value = 1.042258577816105
$ws.Range("K11").Value = 1.04272778030424
$ws.Range("L11").Value = 1.038286624803853
$ws.Range("M11").Value = 1.035405011740867

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03536245013538
$ws.Range("D12").Value = 1.038913119561226
$ws.Range("E12").Value = 1.034428788770396
$ws.Range("F12").Value = 1.031332670528146
$ws.Range("I12").Value = 1.037847483211903
$ws.Range("J12").Value = 1.042013228225197
$ws.Range("K12").Value = 1.042514286731846
$ws.Range("L12").Value = 1.038046779540654
$ws.Range("M12").Value = 1.034962371426022

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035433843986003
$ws.Range("D13").Value = 1.03896874849119
$ws.Range("E13").Value = 1.034490139262547
$ws.Range("F13").Value = 1.031437728454885
$ws.Range("I13").Value = 1.037869987841229
$ws.Range("J13").Value = 1.042065874642312
$ws.Range("K13").Value = 1.042560099640274
$ws.Range("L13").Value = 1.038098242658869
$ws.Range("M13").Value = 1.035057341753489

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035667714398853
$ws.Range("D14").Value = 1.039150971449582
$ws.Range("E14").Value = 1.034691117285293
$ws.Range("F14").Value = 1.031781883203502
$ws.Range("I14").Value = 1.037943660223457
$ws.Range("J14").Value = 1.042238305291427
$ws.Range("K14").Value = 1.042710140881155
$ws.Range("L14").Value = 1.038266806096487
$ws.Range("M14").Value = 1.035368433063776

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035811766323142
$ws.Range("D15").Value = 1.039263207348343
$ws.Range("E15").Value = 1.03481491488677
$ws.Range("F15").Value = 1.03199387095742
$ws.Range("I15").Value = 1.03798900197167
$ws.Range("J15").Value = 1.042344492619781
$ws.Range("K15").Value = 1.042802534172369
$ws.Range("L15").Value = 1.038370618317684
$ws.Range("M15").Value = 1.03556004114329

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036649554884634
$ws.Range("D16").Value = 1.039915897569308
$ws.Range("E16").Value = 1.035534992608689
$ws.Range("F16").Value = 1.033226869097917
$ws.Range("I16").Value = 1.038252145671138
$ws.Range("J16").Value = 1.042961745961935
$ws.Range("K16").Value = 1.043339511968103
$ws.Range("L16").Value = 1.038974163276145
$ws.Range("M16").Value = 1.036674294283109

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.037174506422346
$ws.Range("D17").Value = 1.040324814245112
$ws.Range("E17").Value = 1.035986262403656
$ws.Range("F17").Value = 1.033999551324974
$ws.Range("I17").Value = 1.038416537485011
$ws.Range("J17").Value = 1.043348230684766
$ws.Range("K17").Value = 1.043675651464849
$ws.Range("L17").Value = 1.039352150597754
$ws.Range("M17").Value = 1.037372378100631

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.037480498552709
$ws.Range("D18").Value = 1.04056315046438
$ws.Range("E18").Value = 1.03624933247421
$ws.Range("F18").Value = 1.034449981021944
$ws.Range("I18").Value = 1.038512183469583
$ws.Range("J18").Value = 1.043573409696201
$ws.Range("K18").Value = 1.04387146806312
$ws.Range("L18").Value = 1.039572409620497
$ws.Range("M18").Value = 1.03777925460121

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.037584799734647
$ws.Range("D19").Value = 1.040644386934441
$ws.Range("E19").Value = 1.03633900765288
$ws.Range("F19").Value = 1.034603521870611
$ws.Range("I19").Value = 1.03854475549217
$ws.Range("J19").Value = 1.043650147534905
$ws.Range("K19").Value = 1.043938194574525
$ws.Range("L19").Value = 1.039647476027717
$ws.Range("M19").Value = 1.037917937894383

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037118205170541
$ws.Range("D20").Value = 1.040280959803585
$ws.Range("E20").Value = 1.03593786077511
$ws.Range("F20").Value = 1.033916677145478
$ws.Range("I20").Value = 1.038398924746609
$ws.Range("J20").Value = 1.043306790538943
$ws.Range("K20").Value = 1.043639612556752
$ws.Range("L20").Value = 1.039311618373345
$ws.Range("M20").Value = 1.03729751187399

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035598854437976
$ws.Range("D21").Value = 1.039097319189624
$ws.Range("E21").Value = 1.034631940858973
$ws.Range("F21").Value = 1.031680550186303
$ws.Range("I21").Value = 1.037921976014176
$ws.Range("J21").Value = 1.042187539791546
$ws.Range("K21").Value = 1.042665968332122
$ws.Range("L21").Value = 1.03821717780776
$ws.Range("M21").Value = 1.035276838109974

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.03464169185667
$ws.Range("D22").Value = 1.038351476570501
$ws.Range("E22").Value = 1.033809482690574
$ws.Range("F22").Value = 1.030272118944649
$ws.Range("I22").Value = 1.037619908275768
$ws.Range("J22").Value = 1.041481518360323
$ws.Range("K22").Value = 1.042051530332815
$ws.Range("L22").Value = 1.037527087317498
$ws.Range("M22").Value = 1.034003511470126

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035149283378646
$ws.Range("D23").Value = 1.038747019212882
$ws.Range("E23").Value = 1.034245615633044
$ws.Range("F23").Value = 1.031018996920251
$ws.Range("I23").Value = 1.0377802487767
$ws.Range("J23").Value = 1.041856014169128
$ws.Range("K23").Value = 1.042377472396081
$ws.Range("L23").Value = 1.037893106195083
$ws.Range("M23").Value = 1.034678801426716

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.037143645910016
$ws.Range("D24").Value = 1.040300776287467
$ws.Range("E24").Value = 1.03595973184435
$ws.Range("F24").Value = 1.033954125230035
$ws.Range("I24").Value = 1.038406883929308
$ws.Range("J24").Value = 1.043325516332057
$ws.Range("K24").Value = 1.043655897753533
$ws.Range("L24").Value = 1.03932993380344
$ws.Range("M24").Value = 1.03733134163685

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039450048590257
$ws.Range("D25").Value = 1.04209686261007
$ws.Range("E25").Value = 1.037943084584006
$ws.Range("F25").Value = 1.037349931550897
$ws.Range("I25").Value = 1.039124616066789
$ws.Range("J25").Value = 1.045020967466597
$ws.Range("K25").Value = 1.04512973325593
$ws.Range("L25").Value = 1.04098889395836
$ws.Range("M25").Value = 1.040397597793815
